$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "VOTING 3 BEST FBETA corrected" -- row 2 holds the voting_3best_fbeta_soft
# results; update the metric values with the corrected figures.
$ws.Range("B2").Value = 0.78010686395396633
$ws.Range("C2").Value = 0.88820961264509801
$ws.Range("D2").Value = 0.87951807228915657
$ws.Range("E2").Value = 0.77659574468085102
$ws.Range("F2").Value = 0.80683498090524919
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 21
$ws.Range("I2").Value = 1083.818585634232
$ws.Range("J2").Value = 8.5441548824310303

# Reflect the corrected row as the active selection, matching the
# author's last interaction with the sheet.
$ws.Range("A2:J2").Select()
